# feat: add 2022-Q1 data
#
# 1) The existing "总计" sheet (sheetId=6) is renamed to "2022-Q1" and its
#    content is replaced with the 2022-Q1 fund-holding detail rows.
# 2) A brand-new "总计" sheet is inserted right after it (sheetId=7) with the
#    historical summary table plus the new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet (last sheet) as the new "2022-Q1"
# detail sheet.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = "2022-Q1"

# Clear any previous contents just in case.
$q1.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$q1rows = @(
    @("519655", "银河现代服务主题灵活配置混合", "3.89", "85.22", "7.40", "0.2879", 1),
    @("014185", "招商专精特新股票A",           "8.37", "30.94", "1.65", "0.1381", 8),
    @("570005", "诺德成长优势混合",             "4.30", "62.89", "2.92", "0.1256", 10),
    @("217013", "招商中小盘精选混合",           "2.53", "83.36", "3.99", "0.1009", 9),
    @("012036", "诺德兴远优选一年持有期混合型证券投资基金", "2.75", "52.19", "2.76", "0.0759", 7),
    @("014186", "招商专精特新股票C",           "3.46", "30.94", "1.65", "0.0571", 8),
    @("003561", "诺德成长精选灵活配置混合A",     "0.60", "53.97", "2.95", "0.0177", 7),
    @("006718", "国融融盛龙头严选混合A",        "0.11", "89.32", "4.47", "0.0049", 9),
    @("003562", "诺德成长精选灵活配置混合C",     "0.00", "53.97", "2.95", 0,        7),
    @("006719", "国融融盛龙头严选混合C",        "0.00", "89.32", "4.47", 0,        9)
)

for ($r = 0; $r -lt $q1rows.Length; $r++) {
    $row = $q1rows[$r]
    $excelRow = $r + 2

    $q1.Cells.Item($excelRow, 1).Value = $r

    $q1.Cells.Item($excelRow, 2).Value = "'" + $row[0]
    $q1.Cells.Item($excelRow, 3).Value = "'" + $row[1]
    $q1.Cells.Item($excelRow, 4).Value = "'" + $row[2]
    $q1.Cells.Item($excelRow, 5).Value = "'" + $row[3]
    $q1.Cells.Item($excelRow, 6).Value = "'" + $row[4]

    $g = $row[5]
    if ($g -eq 0) {
        $q1.Cells.Item($excelRow, 7).Value = 0
    } else {
        $q1.Cells.Item($excelRow, 7).Value = "'" + $g
    }

    $q1.Cells.Item($excelRow, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Step 2: add a brand-new "总计" sheet right after "2022-Q1" with the
# updated summary table (old rows + the new 2022-Q1 row on top).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 10, 0.81),
    @("2021-Q4", 5, 1.36),
    @("2021-Q3", 3, 1.53),
    @("2021-Q2", 7, 5.24),
    @("2021-Q1", 7, 2.98),
    @("2020-Q4", 1, 1.06)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = $totalRows[$r]
    $excelRow = $r + 2

    $total.Cells.Item($excelRow, 1).Value = $r
    $total.Cells.Item($excelRow, 2).Value = "'" + $row[0]
    $total.Cells.Item($excelRow, 3).Value = $row[1]
    $total.Cells.Item($excelRow, 4).Value = $row[2]
}
